$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at row 21 (pushes the existing "07.22.19" block from
# rows 21-46 down to rows 28-53, copying row 20's formatting into the
# freshly inserted rows, matching Excel's native Insert behaviour).
$ws.Rows("21:27").Insert()

# The new rows describe 7 additional "08.09.19" s1cDNA samples (#20-#26)
# that were appended to the first sequencing run block.
$newRows = @(
    @{C=20; J=437; Fmt='plain';  K="Brent_2c-20_GTAC_60_SIC_Index2_6_CTGCAAT_GACCTTGT_S21_R1_001.fastq.gz"},
    @{C=21; J=0;   Fmt='numfmt'; K="Brent_2c-21_GTAC_61_SIC_Index2_6_CAAGCCG_GACCTTGT_S22_R1_001.fastq.gz"},
    @{C=22; J=2;   Fmt='plain';  K="Brent_2c-22_GTAC_62_SIC_Index2_6_GGGTCAA_GACCTTGT_S23_R1_001.fastq.gz"},
    @{C=23; J=0;   Fmt='numfmt'; K="Brent_2c-23_GTAC_63_SIC_Index2_6_GCAACGC_GACCTTGT_S24_R1_001.fastq.gz"},
    @{C=24; J=1;   Fmt='plain';  K="Brent_2c-24_GTAC_64_SIC_Index2_6_TGATTAC_GACCTTGT_S25_R1_001.fastq.gz"},
    @{C=25; J=2;   Fmt='plain';  K="Brent_2c-25_GTAC_65_SIC_Index2_6_TGCTGGG_GACCTTGT_S26_R1_001.fastq.gz"},
    @{C=26; J=9;   Fmt='plain';  K="Brent_2c-26_GTAC_66_SIC_Index2_6_GACACAG_GACCTTGT_S27_R1_001.fastq.gz"}
)

# Column A holds a "MM.DD.YY"-shaped label that Excel's auto-detection
# would otherwise coerce into a date serial. Force it to text first...
$colA = $ws.Range("A21:A27")
$colA.NumberFormat = "@"

$r = 21
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = "08.09.19"
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = 3906
    $ws.Cells.Item($r, 5).Value = "NextSeq"
    $ws.Cells.Item($r, 6).Value = "MidOutput"
    $ws.Cells.Item($r, 7).Value = "fullRNASeq"
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 2

    $jc = $ws.Cells.Item($r, 10)
    $jc.Value = $row.J
    $jc.Font.Name = "Helvetica"
    $jc.Font.Size = 12
    $jc.Font.Color = 0
    if ($row.Fmt -eq 'numfmt') {
        $jc.NumberFormat = "#,##0"
    }

    $ws.Cells.Item($r, 11).Value = $row.K

    $r = $r + 1
}

# ...then restore column A's original (General/date-neutral) formatting by
# pasting row 20's format over it, so it ends up styled exactly like the
# rest of the libraryDate column instead of keeping a bespoke "@" style.
$ws.Cells.Item(20, 1).Copy()
$colA.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where Excel would naturally land after this edit.
$ws.Cells.Item(28, 9).Select()
